$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the header cell first, while it is still the sole user of its style,
# so the shared style/font gets reused instead of forking.
$ws.Range("A1").Font.Bold = $true

# Copy the (now bold) header formatting from A1 into the new B1 header cell.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rename existing header and add the Spanish header
$ws.Range("A1").Value = "EN_URL"
$ws.Range("B1").Value = "ES_URL"

# Existing English URLs stay in column A
$ws.Range("A2").Value = "https://www.rybelsus.com/"
$ws.Range("A3").Value = "https://www.rybelsus.com/savings-and-support.html"

# New Spanish URLs go in column B
$ws.Range("B2").Value = "https://espanol.rybelsus.com/"
$ws.Range("B3").Value = "https://espanol.rybelsus.com/ahorros-y-apoyo.html"

# Match the selection left by the edit
$ws.Range("B1").Select()
